# Update recalculated NATMI ligand-receptor edge-weight statistics
# following a re-run of the pipeline with refreshed TPM expression values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("G2").Value = 0.02354566666666667
$ws.Range("H2").Value = 0.07063700000000001
$ws.Range("I2").Value = 0.002815555392485919
$ws.Range("J2").Value = 0.002815555392485918
$ws.Range("M2").Value = 1.01111
$ws.Range("N2").Value = 3.03333
$ws.Range("O2").Value = 0.04063212692754557
$ws.Range("P2").Value = 0.04063212692754556
$ws.Range("Q2").Value = 0.02380725902333334
$ws.Range("R2").Value = 0.21426533121
$ws.Range("S2").Value = 0.0001144020040790232
$ws.Range("T2").Value = 0.0001144020040790232

# Row 3: ECs -> FAPs
$ws.Range("G3").Value = 0.02354566666666667
$ws.Range("H3").Value = 0.07063700000000001
$ws.Range("I3").Value = 0.002815555392485919
$ws.Range("J3").Value = 0.002815555392485918
$ws.Range("O3").Value = 0.4065982422683317
$ws.Range("P3").Value = 0.4065982422683317
$ws.Range("Q3").Value = 0.2382348748165556
$ws.Range("R3").Value = 2.144113873349
$ws.Range("S3").Value = 0.001144799873593897
$ws.Range("T3").Value = 0.001144799873593897

# Row 4: ECs -> MuSCs
$ws.Range("G4").Value = 0.02354566666666667
$ws.Range("H4").Value = 0.07063700000000001
$ws.Range("I4").Value = 0.002815555392485919
$ws.Range("J4").Value = 0.002815555392485918
$ws.Range("O4").Value = 0.5527696308041227
$ws.Range("P4").Value = 0.5527696308041226
$ws.Range("Q4").Value = 0.3238799141441112
$ws.Range("R4").Value = 2.914919227297001
$ws.Range("S4").Value = 0.001556353514812998
$ws.Range("T4").Value = 0.001556353514812997

# Row 5: FAPs -> ECs
$ws.Range("I5").Value = 0.9868456480383168
$ws.Range("J5").Value = 0.9868456480383166
$ws.Range("M5").Value = 1.01111
$ws.Range("N5").Value = 3.03333
$ws.Range("O5").Value = 0.04063212692754557
$ws.Range("P5").Value = 0.04063212692754556
$ws.Range("Q5").Value = 8.344389182183333
$ws.Range("R5").Value = 75.09950263965
$ws.Range("S5").Value = 0.04009763762898885
$ws.Range("T5").Value = 0.04009763762898883

# Row 6: FAPs -> FAPs
$ws.Range("I6").Value = 0.9868456480383168
$ws.Range("J6").Value = 0.9868456480383166
$ws.Range("O6").Value = 0.4065982422683317
$ws.Range("P6").Value = 0.4065982422683317
$ws.Range("S6").Value = 0.4012497058825323
$ws.Range("T6").Value = 0.4012497058825323

# Row 7: FAPs -> MuSCs
$ws.Range("I7").Value = 0.9868456480383168
$ws.Range("J7").Value = 0.9868456480383166
$ws.Range("O7").Value = 0.5527696308041227
$ws.Range("P7").Value = 0.5527696308041226
$ws.Range("S7").Value = 0.5454983045267955
$ws.Range("T7").Value = 0.5454983045267954

# Row 8: MuSCs -> ECs
$ws.Range("I8").Value = 0.0103387965691973
$ws.Range("J8").Value = 0.0103387965691973
$ws.Range("M8").Value = 1.01111
$ws.Range("N8").Value = 3.03333
$ws.Range("O8").Value = 0.04063212692754557
$ws.Range("P8").Value = 0.04063212692754556
$ws.Range("Q8").Value = 0.08742090763666664
$ws.Range("R8").Value = 0.7867881687299999
$ws.Range("S8").Value = 0.0004200872944776975
$ws.Range("T8").Value = 0.0004200872944776974

# Row 9: MuSCs -> FAPs
$ws.Range("I9").Value = 0.0103387965691973
$ws.Range("J9").Value = 0.0103387965691973
$ws.Range("O9").Value = 0.4065982422683317
$ws.Range("P9").Value = 0.4065982422683317
$ws.Range("S9").Value = 0.004203736512205482
$ws.Range("T9").Value = 0.004203736512205482

# Row 10: MuSCs -> MuSCs
$ws.Range("I10").Value = 0.0103387965691973
$ws.Range("J10").Value = 0.0103387965691973
$ws.Range("O10").Value = 0.5527696308041227
$ws.Range("P10").Value = 0.5527696308041226
$ws.Range("S10").Value = 0.005714972762514123
$ws.Range("T10").Value = 0.005714972762514122
